# Apply the bill update to Sheet2: add new rows recording the RMB 5620
# payment and related entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate() | Out-Null

$ws.Range("A27").Value = "RMB 5620"

$ws.Range("A28").Value = 5620
$ws.Range("B28").Value = "RMB"

$ws.Range("A30").Value = 190
$ws.Range("B30").Value = "RMB"

$ws.Range("E33").Value = 1060
$ws.Range("F33").Value = "RMB"

$ws.Range("A34").Value = 5810
$ws.Range("B34").Value = "RMB"

$ws.Range("A36").Value = 6870
$ws.Range("B36").Value = "RMB"

$ws.Range("C36").Select() | Out-Null
